$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dependencies")
$r = $ws.Range("B6")
for ($i=5; $i -le 12; $i++) {
  try {
    $b = $r.Borders.Item($i)
    Write-Host "$i LineStyle=$($b.LineStyle) Weight=$($b.Weight) Color=$($b.Color) ColorIndex=$($b.ColorIndex)"
  } catch {
    Write-Host "$i ERROR $_"
  }
}
